$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 10883.167
$ws.Range("I33").Value = 11517.679
$ws.Range("J33").Value = 2000
$ws.Range("K33").Value = 11517.679
$ws.Range("L33").Value = 2000
$ws.Range("M33").Value = -11288.679
$ws.Range("N33").Value = -2458
$ws.Range("H86").Value = 999
$ws.Range("I86").Value = 999
$ws.Range("K86").Value = 999
$ws.Range("M86").Value = 124
$ws.Range("H89").Value = 999
$ws.Range("I89").Value = 999
$ws.Range("K89").Value = 4995
$ws.Range("M89").Value = 621
$ws.Range("H98").Value = 2425.1516
$ws.Range("I98").Value = 2493.125
$ws.Range("K98").Value = 2493.125
$ws.Range("M98").Value = -995.125
$ws.Range("H122").Value = 2425.1516
$ws.Range("I122").Value = 2493.125
$ws.Range("K122").Value = 7479.375
$ws.Range("M122").Value = -5029.375
$ws.Range("H132").Value = 1542.3158
$ws.Range("I132").Value = 1260.6
$ws.Range("K132").Value = 3781.8
$ws.Range("M132").Value = -1251.8
$ws.Range("H137").Value = 1922.6666
$ws.Range("I137").Value = 1922.6666
$ws.Range("K137").Value = 5767.9998
$ws.Range("M137").Value = -3217.9998
$ws.Range("H138").Value = 8273.5
$ws.Range("I138").Value = 8273.5
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 24820.5
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -19680.5
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3064.9219
$ws.Range("I32").Value = 2557.5356
$ws.Range("J32").Value = 6616.625
$ws.Range("K32").Value = 2557.5356
$ws.Range("L32").Value = 6616.625
$ws.Range("M32").Value = -2270.5356
$ws.Range("N32").Value = -7190.625
$ws.Range("H132").Value = 4121.8335
$ws.Range("I132").Value = 3315.7812
$ws.Range("K132").Value = 9947.3436
$ws.Range("M132").Value = -7417.3436

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3240.2144
$ws.Range("I20").Value = 3124
$ws.Range("J20").Value = 3666.3333
$ws.Range("K20").Value = 3124
$ws.Range("L20").Value = 3666.3333
$ws.Range("M20").Value = -2877
$ws.Range("N20").Value = -4160.3333
$ws.Range("H35").Value = 33940
$ws.Range("J35").Value = 33940
$ws.Range("L35").Value = 33940
$ws.Range("N35").Value = -34560
$ws.Range("H105").Value = 8199.299999999999
$ws.Range("I105").Value = 7416.5
$ws.Range("J105").Value = 9373.5
$ws.Range("K105").Value = 7416.5
$ws.Range("L105").Value = 9373.5
$ws.Range("M105").Value = -5669.5
$ws.Range("N105").Value = -12867.5
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 6065.2
$ws.Range("I134").Value = 5976.5264
$ws.Range("J134").Value = 7750
$ws.Range("K134").Value = 17929.5792
$ws.Range("L134").Value = 23250
$ws.Range("M134").Value = -15394.5792
$ws.Range("N134").Value = -28320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4912.5
$ws.Range("I31").Value = 3849.2856
$ws.Range("J31").Value = 5975.7144
$ws.Range("K31").Value = 3849.2856
$ws.Range("L31").Value = 5975.7144
$ws.Range("M31").Value = -3554.2856
$ws.Range("N31").Value = -6565.7144
$ws.Range("H34").Value = 4912.5
$ws.Range("I34").Value = 3849.2856
$ws.Range("J34").Value = 5975.7144
$ws.Range("K34").Value = 3849.2856
$ws.Range("L34").Value = 5975.7144
$ws.Range("M34").Value = -3647.2856
$ws.Range("N34").Value = -6379.7144
$ws.Range("H58").Value = 8838.888999999999
$ws.Range("I58").Value = 5400
$ws.Range("J58").Value = 10558.333
$ws.Range("K58").Value = 5400
$ws.Range("L58").Value = 10558.333
$ws.Range("M58").Value = -5197
$ws.Range("N58").Value = -10964.333
$ws.Range("H86").Value = 5757.5
$ws.Range("I86").Value = 5134.4
$ws.Range("J86").Value = 7092.7144
$ws.Range("K86").Value = 5134.4
$ws.Range("L86").Value = 7092.7144
$ws.Range("M86").Value = -4011.4
$ws.Range("N86").Value = -9338.714400000001
$ws.Range("H89").Value = 5757.5
$ws.Range("I89").Value = 5134.4
$ws.Range("J89").Value = 7092.7144
$ws.Range("K89").Value = 25672
$ws.Range("L89").Value = 35463.572
$ws.Range("M89").Value = -20056
$ws.Range("N89").Value = -46695.572
$ws.Range("H107").Value = 2005.5555
$ws.Range("I107").Value = 1939.2667
$ws.Range("J107").Value = 2088.4167
$ws.Range("K107").Value = 1939.2667
$ws.Range("L107").Value = 2088.4167
$ws.Range("M107").Value = -19.2666999999999
$ws.Range("N107").Value = -5928.4167
$ws.Range("H134").Value = 9655.368
$ws.Range("I134").Value = 7167
$ws.Range("J134").Value = 10318.934
$ws.Range("K134").Value = 21501
$ws.Range("L134").Value = 30956.802
$ws.Range("M134").Value = -18966
$ws.Range("N134").Value = -36026.802
$ws.Range("H136").Value = 8838.888999999999
$ws.Range("I136").Value = 5400
$ws.Range("J136").Value = 10558.333
$ws.Range("K136").Value = 16200
$ws.Range("L136").Value = 31674.999
$ws.Range("M136").Value = -13650
$ws.Range("N136").Value = -36774.999
$ws.Range("H141").Value = 39231.562
$ws.Range("J141").Value = 39231.562
$ws.Range("L141").Value = 39231.562
$ws.Range("N141").Value = -49591.562

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1999.5
$ws.Range("J75").Value = 1999
$ws.Range("L75").Value = 5997
$ws.Range("N75").Value = -7993
$ws.Range("H78").Value = 1999.5
$ws.Range("J78").Value = 1999
$ws.Range("L78").Value = 17991
$ws.Range("N78").Value = -27975
$ws.Range("H122").Value = 1578.4138
$ws.Range("I122").Value = 1279.8
$ws.Range("J122").Value = 1735.579
$ws.Range("K122").Value = 11518.2
$ws.Range("L122").Value = 15620.211
$ws.Range("M122").Value = -9068.199999999999
$ws.Range("N122").Value = -20520.211

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4963.3335
$ws.Range("I132").Value = 4945
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 14835
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -12305
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5065.6665
$ws.Range("I40").Value = 5331.3335
$ws.Range("K40").Value = 5331.3335
$ws.Range("M40").Value = -5195.3335
$ws.Range("H43").Value = 358333.34
$ws.Range("J43").Value = 358333.34
$ws.Range("L43").Value = 358333.34
$ws.Range("N43").Value = -358719.34
$ws.Range("H122").Value = 4442.5557
$ws.Range("I122").Value = 4617.1665
$ws.Range("K122").Value = 13851.4995
$ws.Range("M122").Value = -11401.4995
$ws.Range("H132").Value = 12350.8
$ws.Range("I132").Value = 21901.6
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 65704.79999999999
$ws.Range("L132").Value = 8400
$ws.Range("M132").Value = -63174.79999999999
$ws.Range("N132").Value = -13460
$ws.Range("H133").Value = 72972.5
$ws.Range("J133").Value = 72972.5
$ws.Range("L133").Value = 72972.5
$ws.Range("N133").Value = -78032.5
$ws.Range("H137").Value = 85100
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 40496.668
$ws.Range("J70").Value = 40496.668
$ws.Range("L70").Value = 40496.668
$ws.Range("N70").Value = -41126.668
$ws.Range("H73").Value = 40496.668
$ws.Range("J73").Value = 40496.668
$ws.Range("L73").Value = 40496.668
$ws.Range("N73").Value = -42680.668
$ws.Range("H122").Value = 4038.2964
$ws.Range("I122").Value = 3524.4614
$ws.Range("J122").Value = 4515.4287
$ws.Range("K122").Value = 10573.3842
$ws.Range("L122").Value = 13546.2861
$ws.Range("M122").Value = -8123.3842
$ws.Range("N122").Value = -18446.2861
$ws.Range("H126").Value = 4872.433
$ws.Range("I126").Value = 2417.8096
$ws.Range("K126").Value = 7253.4288
$ws.Range("M126").Value = -4783.4288
$ws.Range("H132").Value = 3442.4285
$ws.Range("I132").Value = 2187.3704
$ws.Range("K132").Value = 6562.111199999999
$ws.Range("M132").Value = -4032.111199999999
